$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header labels (shared-string reindex only; word changes) ---
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"

# --- Rows 3-7: row labels shift down by one (new "fraud" row inserted) ---
$ws.Range("A4").Value = "fraud"
$ws.Range("A5").Value = "crisis"
$ws.Range("A6").Value = "panic"
$ws.Range("A7").Value = "sc"

# --- Rows 3-7: update B/C/D/H metrics (E,F,G unchanged) ---
$ws.Range("B3").Value = 0.8529411764705882
$ws.Range("C3").Value = 29
$ws.Range("D3").Value = 29
$ws.Range("H3").Value = 5
$ws.Range("B4").Value = 0.6944444444444444
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("H4").Value = 11
$ws.Range("B5").Value = 0.6232876712328768
$ws.Range("C5").Value = 182
$ws.Range("D5").Value = 182
$ws.Range("H5").Value = 110
$ws.Range("B6").Value = 0.2131782945736434
$ws.Range("C6").Value = 110
$ws.Range("D6").Value = 110
$ws.Range("H6").Value = 406
$ws.Range("B7").Value = 0.1957671957671958
$ws.Range("C7").Value = 37
$ws.Range("D7").Value = 37
$ws.Range("H7").Value = 152

# --- Clear old row 8 left-table entries (A8:H8) - no longer present ---
$ws.Range("A8:H8").Clear()

# --- Right-hand table J3:Q29: word list grew by one row (insert before, reflow down) ---
$ws.Range("J3").Value = "interesting"
$ws.Range("K3").Value = 0.9393939393939394
$ws.Range("L3").Value = 31
$ws.Range("M3").Value = 31
$ws.Range("Q3").Value = 2
$ws.Range("J4").Value = "best"
$ws.Range("K4").Value = 0.9322033898305084
$ws.Range("L4").Value = 55
$ws.Range("M4").Value = 55
$ws.Range("Q4").Value = 4
$ws.Range("J5").Value = "great"
$ws.Range("K5").Value = 0.875
$ws.Range("L5").Value = 98
$ws.Range("M5").Value = 98
$ws.Range("Q5").Value = 14
$ws.Range("J6").Value = "positive"
$ws.Range("K6").Value = 0.8448275862068966
$ws.Range("L6").Value = 49
$ws.Range("M6").Value = 49
$ws.Range("Q6").Value = 9
$ws.Range("J7").Value = "love"
$ws.Range("K7").Value = 0.8260869565217391
$ws.Range("L7").Value = 38
$ws.Range("M7").Value = 38
$ws.Range("Q7").Value = 8
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8170731707317073
$ws.Range("L8").Value = 67
$ws.Range("M8").Value = 67
$ws.Range("Q8").Value = 15
$ws.Range("J9").Value = "thank"
$ws.Range("K9").Value = 0.8125
$ws.Range("L9").Value = 104
$ws.Range("M9").Value = 104
$ws.Range("Q9").Value = 24
$ws.Range("J10").Value = "special"
$ws.Range("K10").Value = 0.8055555555555556
$ws.Range("L10").Value = 29
$ws.Range("M10").Value = 29
$ws.Range("Q10").Value = 7
$ws.Range("J11").Value = "free"
$ws.Range("K11").Value = 0.8
$ws.Range("L11").Value = 96
$ws.Range("M11").Value = 96
$ws.Range("Q11").Value = 24
$ws.Range("J12").Value = "safe"
$ws.Range("K12").Value = 0.7535211267605634
$ws.Range("L12").Value = 107
$ws.Range("M12").Value = 107
$ws.Range("Q12").Value = 35
$ws.Range("J13").Value = "support"
$ws.Range("K13").Value = 0.7075471698113207
$ws.Range("L13").Value = 75
$ws.Range("M13").Value = 75
$ws.Range("Q13").Value = 31
$ws.Range("J14").Value = "confidence"
$ws.Range("K14").Value = 0.6944444444444444
$ws.Range("L14").Value = 25
$ws.Range("M14").Value = 25
$ws.Range("Q14").Value = 11
$ws.Range("J15").Value = "good"
$ws.Range("K15").Value = 0.675
$ws.Range("L15").Value = 108
$ws.Range("M15").Value = 108
$ws.Range("Q15").Value = 52
$ws.Range("J16").Value = "safety"
$ws.Range("K16").Value = 0.6666666666666666
$ws.Range("L16").Value = 34
$ws.Range("M16").Value = 34
$ws.Range("Q16").Value = 17
$ws.Range("J17").Value = "heroes"
$ws.Range("K17").Value = 0.6595744680851063
$ws.Range("L17").Value = 31
$ws.Range("M17").Value = 31
$ws.Range("Q17").Value = 16
$ws.Range("J18").Value = "relief"
$ws.Range("K18").Value = 0.64
$ws.Range("L18").Value = 32
$ws.Range("M18").Value = 32
$ws.Range("Q18").Value = 18
$ws.Range("J19").Value = "well"
$ws.Range("K19").Value = 0.6382978723404256
$ws.Range("L19").Value = 60
$ws.Range("M19").Value = 60
$ws.Range("Q19").Value = 34
$ws.Range("J20").Value = "fresh"
$ws.Range("K20").Value = 0.6041666666666666
$ws.Range("L20").Value = 29
$ws.Range("M20").Value = 29
$ws.Range("Q20").Value = 19
$ws.Range("J21").Value = "better"
$ws.Range("K21").Value = 0.6031746031746031
$ws.Range("L21").Value = 38
$ws.Range("M21").Value = 38
$ws.Range("Q21").Value = 25
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5535248041775457
$ws.Range("L22").Value = 212
$ws.Range("M22").Value = 212
$ws.Range("Q22").Value = 171
$ws.Range("J23").Value = "like"
$ws.Range("K23").Value = 0.4764705882352941
$ws.Range("L23").Value = 162
$ws.Range("M23").Value = 162
$ws.Range("Q23").Value = 178
$ws.Range("J24").Value = "care"
$ws.Range("K24").Value = 0.4606741573033708
$ws.Range("L24").Value = 41
$ws.Range("M24").Value = 41
$ws.Range("Q24").Value = 48
$ws.Range("J25").Value = "increase"
$ws.Range("K25").Value = 0.4487179487179487
$ws.Range("L25").Value = 35
$ws.Range("M25").Value = 35
$ws.Range("Q25").Value = 43
$ws.Range("J26").Value = "help"
$ws.Range("K26").Value = 0.4440677966101695
$ws.Range("L26").Value = 131
$ws.Range("M26").Value = 131
$ws.Range("Q26").Value = 164
$ws.Range("J27").Value = "protect"
$ws.Range("K27").Value = 0.3972602739726027
$ws.Range("L27").Value = 29
$ws.Range("M27").Value = 29
$ws.Range("Q27").Value = 44
$ws.Range("J28").Value = "please"
$ws.Range("K28").Value = 0.3389121338912134
$ws.Range("L28").Value = 81
$ws.Range("M28").Value = 81
$ws.Range("Q28").Value = 158
$ws.Range("J29").Value = "store"
$ws.Range("K29").Value = 0.02796420581655481
$ws.Range("L29").Value = 25
$ws.Range("M29").Value = 25
$ws.Range("Q29").Value = 869

# --- New row 29 needs N/O/P plus style copied from row 28 (border+bold+center) ---
$ws.Range("N29").Value = 1
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = $false

# Copy formatting (border/bold/alignment) from J28:Q28 onto J29:Q29 so style index matches
$ws.Range("J28:Q28").Copy() | Out-Null
$ws.Range("J29:Q29").PasteSpecial(-4122) | Out-Null

Write-Output "edit applied"
